$d = $word.ActiveDocument

# Remove the _GoBack bookmark left over from the previous edit session.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Collapse "Van't " + "Hoff's " (with the spell-check proofing marks around
# "Van't") into a single run / remove the proofing-error markers, matching
# how Word normalises the paragraph when it is retyped/resaved.
$d.Content.Find.Execute("Van") | Out-Null
